# Trade #78 closed at 2026-02-18 00:28:15 - unknown UNKNOWN +0.000%
#
# 1) Trade #106 (row 107 in "All Trades", row 25 in "momentum") closes:
#    OPEN -> CLOSED, exit price 0.95, small positive P&L, exit reason "early_exit".
# 2) A brand-new Trade #135 (MarketMaking, still OPEN) gets appended to
#    "All Trades" (row 136) and to "MarketMaking" (row 56).
# 3) The "Summary" and "Strategy Status" roll-up sheets are refreshed to
#    reflect the above.

$wb = $excel.ActiveWorkbook

# Excel auto-detects strings that look like dates (e.g. "2026-02-18") and
# silently converts them to date serials. Forcing a text number format
# first keeps them as literal strings, matching the source data (every
# other Date/Time/Reason cell in these sheets is plain text).
function Set-Text($ws, $row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
}

function Set-Plain($ws, $row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = $text
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.46
$summary.Range("B4").Value = 0.57
$summary.Range("B6").Value = 106
$summary.Range("B7").Value = 52
$summary.Range("B9").Value = 49.06

# ---------------------------------------------------------------------
# Strategy Status sheet - "momentum" row (row 11)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C11").Value = 99.3
$status.Range("D11").Value = 24
$status.Range("E11").Value = -0.7
$status.Range("F11").Value = -0.7
$status.Range("G11").Value = 29.17

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close out trade #106 (row 107): G=Exit Price, H=Status, I=P&L%, J=P&L$,
# K=Capital After, L=Exit Reason, M=Duration (min)
$allTrades.Cells.Item(107, 7).Value = 0.95
Set-Plain $allTrades 107 8 "CLOSED"
$allTrades.Cells.Item(107, 9).Value = 1.0638
$allTrades.Cells.Item(107, 10).Value = 0.01
$allTrades.Cells.Item(107, 11).Value = 99.3
Set-Plain $allTrades 107 12 "early_exit"
$allTrades.Cells.Item(107, 13).Value = 0.13

# Append new trade #135 (row 136)
$allTrades.Cells.Item(136, 1).Value = 135
Set-Text $allTrades 136 2 "2026-02-18"
Set-Plain $allTrades 136 3 "00:28:09"
Set-Plain $allTrades 136 4 "MarketMaking"
Set-Plain $allTrades 136 5 "DOWN"
$allTrades.Cells.Item(136, 6).Value = 0.9399999999999999
Set-Plain $allTrades 136 8 "OPEN"
$allTrades.Cells.Item(136, 9).Value = 0
$allTrades.Cells.Item(136, 10).Value = 0
$allTrades.Cells.Item(136, 11).Value = 99.47967800952271
$allTrades.Cells.Item(136, 13).Value = 0
$allTrades.Cells.Item(136, 14).Value = 0
$allTrades.Cells.Item(136, 15).Value = 0
$allTrades.Cells.Item(136, 16).Value = 0.65
Set-Plain $allTrades 136 17 "Wide spread capture: 392 bps vs avg 311 bps"

# ---------------------------------------------------------------------
# momentum sheet - trade #106 (row 25)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Cells.Item(25, 7).Value = 0.95
Set-Plain $momentum 25 8 "CLOSED"
$momentum.Cells.Item(25, 9).Value = 1.0638
$momentum.Cells.Item(25, 10).Value = 0.01
$momentum.Cells.Item(25, 11).Value = 99.3
Set-Plain $momentum 25 16 "early_exit"
$momentum.Cells.Item(25, 17).Value = 0.13

# ---------------------------------------------------------------------
# MarketMaking sheet - append new trade #135 (row 56)
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Cells.Item(56, 1).Value = 135
Set-Text $marketMaking 56 2 "2026-02-18"
Set-Plain $marketMaking 56 3 "00:28:09"
Set-Plain $marketMaking 56 4 "MarketMaking"
Set-Plain $marketMaking 56 5 "DOWN"
$marketMaking.Cells.Item(56, 6).Value = 0.9399999999999999
Set-Plain $marketMaking 56 8 "OPEN"
$marketMaking.Cells.Item(56, 9).Value = 0
$marketMaking.Cells.Item(56, 10).Value = 0
$marketMaking.Cells.Item(56, 11).Value = 99.47967800952271
$marketMaking.Cells.Item(56, 12).Value = 0
$marketMaking.Cells.Item(56, 13).Value = 0
$marketMaking.Cells.Item(56, 14).Value = 0.65
Set-Plain $marketMaking 56 15 "Wide spread capture: 392 bps vs avg 311 bps"
$marketMaking.Cells.Item(56, 17).Value = 0
